$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2,38426,55561779),
    @(3,92274,135259785),
    @(4,31555,46731982),
    @(5,8847,13149356),
    @(6,2039,3030471),
    @(7,158,232093),
    @(12,41875,56801074),
    @(13,9821,14206278),
    @(14,26285,38541108),
    @(15,8405,12472978),
    @(16,2187,3250430),
    @(20,10346,13687389),
    @(21,13588,19611523),
    @(22,32029,47002647),
    @(23,10323,15345349),
    @(24,2672,3972771),
    @(25,521,775592),
    @(27,11843,15811990),
    @(28,7801,11292894),
    @(29,22835,33517093),
    @(30,7895,11744391),
    @(31,1995,2976919),
    @(34,8419,11122360),
    @(35,3317,4791055),
    @(36,7971,11641457),
    @(37,3218,4769961),
    @(39,169,251186),
    @(41,2514,3396922),
    @(42,17536,25355649),
    @(43,51846,75997847),
    @(44,19192,28504392),
    @(45,5680,8455260),
    @(50,16977,22576849),
    @(51,2095,3038378),
    @(52,7121,10464264),
    @(53,2403,3588964),
    @(54,765,1142805),
    @(55,190,281226),
    @(57,7210,9914623),
    @(58,1093,1787129),
    @(59,2711,4437358),
    @(60,1066,1743838),
    @(61,364,599383),
    @(64,1599,2434604),
    @(65,15647,22600505),
    @(66,45347,66352582),
    @(67,15884,23601170),
    @(68,4621,6882551),
    @(69,946,1407168),
    @(73,15289,20142350),
    @(74,53369,77668320),
    @(75,150209,221295399),
    @(76,64938,96765458),
    @(77,20777,31045822),
    @(78,4947,7389401),
    @(85,52667,71593870),
    @(86,4718,6837267),
    @(87,11775,17296986),
    @(88,3933,5862083),
    @(89,1362,2035489),
    @(93,5504,7398810),
    @(94,1631,2350033),
    @(95,5288,7789801),
    @(96,1967,2928826),
    @(98,195,293113),
    @(101,3639,4815769),
    @(102,686,1114825),
    @(103,411,682097),
    @(104,148,242020),
    @(107,10979,15929218),
    @(108,29580,43449509),
    @(109,9897,14714705),
    @(110,2724,4061080),
    @(111,501,746546),
    @(114,9931,13115412),
    @(115,30989,44682006),
    @(116,67000,98042862),
    @(117,21602,32104170),
    @(118,6130,9133021),
    @(124,26187,34958233),
    @(125,36652,52891942),
    @(126,77882,113876117),
    @(127,24118,35797287),
    @(128,6464,9605858),
    @(129,1266,1882311),
    @(133,32230,42785570),
    @(134,13499,19539536),
    @(135,32751,48096690),
    @(136,11602,17237087),
    @(137,2997,4466741),
    @(141,10942,14588547),
    @(142,35771,51664047),
    @(143,82564,120962118),
    @(144,24682,36669243),
    @(145,6481,9671067),
    @(146,1467,2182730),
    @(149,29638,39959063)

)

foreach ($row in $changes) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
